$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2:G13 from "proton" to "p"
$ws.Range("G2:G13").Value = "p"

# Bold + center the header row (A1:K1)
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Select G2:G13 with active cell G2 (mirrors the saved selection state in the file)
$ws.Range("G2:G13").Select()
